$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "loss_to_fup_prob" row (old row 2) entirely - it no longer
# exists in the new layout. Deleting row 2 shifts everything up by one,
# which lines up the remaining rows (lambda, prop_ai, alpha, beta) with
# their new positions (2, 6, 7, 8) before the new rows are inserted. This
# also drops the custom row height / wrap-text formatting that row used to
# carry, reverting row 2 to the sheet's default height.
$ws.Rows.Item(2).Delete()

# Current remaining rows after the delete:
#   row2 = lambda / 0.005 / Boily
#   row3 = prop_ai / 0.15
#   row4 = alpha / 0.9
#   row5 = beta / 3
# Update lambda's value to the new figure.
$ws.Range("B2").Value = 0.003

# Insert three new rows (3, 4, 5) for the new parameters, pushing
# prop_ai/alpha/beta down to rows 6, 7, 8.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Fill in row 5 (cond_rr) first, then rows 3-4, matching the order the
# new labels/notes were authored (preserves shared-string append order).
$ws.Range("A5").Value = "cond_rr"
$ws.Range("B5").Value = 1.8
$ws.Range("D5").Value = "Relative increase in probability of having an HIV-positive male partner among those who report using condoms in previous week at baseline, relative to those who report not using condoms in previous 7 days at baseline."

$ws.Range("A3").Value = "pre_adh_int_rr_bl"
$ws.Range("B3").Value = 0.9

$ws.Range("A4").Value = "pre_adh_int_rr_fu"
$ws.Range("B4").Value = 0.95

# Column A width adjustment (new, longer parameter names no longer use
# auto "best fit" - a fixed width of 16 is set explicitly instead).
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666

# Update selection to match the new active cell recorded in the workbook
$ws.Range("B5").Select()
